# The BLS "unemployment rate" series had no reliable data before 2003, so
# the commit trims the three oldest years (2000-2002) from the data table,
# shifting 2003-2020 up to occupy rows 2-19.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("BLS Data Series")

# Rows 2-4 hold 2000, 2001 and 2002 — delete them outright so every row
# below shifts up three places (2003 becomes row 2, ..., 2020 becomes row 19).
$ws.Range("A2:XFD4").EntireRow.Delete()

# Leave the same 3-row band selected (it now spans the new 2003-2005 rows),
# matching how Excel keeps the selection anchored after a row delete.
$ws.Range("A2:XFD4").Select()
